$d = $word.ActiveDocument

# Locate the heading paragraph "PlantCover_AnalysisDataPrep.Rmd"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "PlantCover_AnalysisDataPrep.Rmd") {
        $target = $p
        break
    }
}

# Insert a brand-new (plain/Normal) paragraph directly after the heading,
# mirroring the description paragraph that follows
# "SpeciesLevel_AnalysisDataPrep.Rmd" above it.
$target.Range.InsertParagraphAfter()
$newPara = $target.Next()
$newPara.Style = "Normal"

$r = $newPara.Range
$r.Collapse(1)
$r.InsertAfter("Compiles climate and")
$r.Collapse(0)
$r.InsertAfter(" total plant cover")
$r.Collapse(0)
$r.InsertAfter(" vegetation data into a table format ready for model input.")
